# Replace use of escape() with encodeURIComponent() in the "queries" sheet
# (auxillaryHash column, rows 2-4), and make "queries" the active sheet/tab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("queries")

$newFormula = "''household_id='+encodeURIComponent(data('household_id'))"

$ws.Range("G2").Value = $newFormula
$ws.Range("G3").Value = $newFormula
$ws.Range("G4").Value = $newFormula

$ws.Activate()
$ws.Select()
$ws.Range("G5").Select()
